$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (New Time) and column E (New Distance count) values
# Row 2 (Total)
$ws.Range("C2").Value = 13.989328384399414
$ws.Range("E2").Value = 2039.0

# Row 4 (Day 0)
$ws.Range("C4").Value = 2.2315924167633057
$ws.Range("E4").Value = 339.0

# Row 5 (Day 1)
$ws.Range("C5").Value = 2.345822334289551
$ws.Range("E5").Value = 343.0

# Row 6 (Day 2)
$ws.Range("C6").Value = 2.624532699584961
$ws.Range("E6").Value = 337.0

# Row 7 (Day 3)
$ws.Range("C7").Value = 2.120365858078003
$ws.Range("E7").Value = 338.0

# Row 8 (Day 4)
$ws.Range("C8").Value = 2.39312481880188
$ws.Range("E8").Value = 339.0

# Row 9 (Day 5)
$ws.Range("C9").Value = 2.2738900184631348
$ws.Range("E9").Value = 343.0
